$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "_old" suffix -> "_FV2304", "_new" suffix -> "_FV2310"
$oldHeaders = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)
$newHeaders = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $oldHeaders[$i]
}

$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeaders[$i]
}

# Turn the data range into an Excel Table (ListObject)
$rng = $ws.Range("A1:U64")
$tbl = $ws.ListObjects.Add(1, $rng, [System.Type]::Missing, 1)
$tbl.Name = "Table1"

# Freeze the header row (split pane above row 2)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
